# Fix typo 'Cyro-EM' to 'Cryo-EM' in the storage_medium lookup sheet.
# The fixed entry keeps its relative ordering right after "Gelatin" (row 11),
# which pushes "DMSO (serum)" and "RNAlater" down by one row each
# (rows 12-14 rotate), matching the author's intended ordering.
# Also bump the pav:createdOn metadata timestamp to reflect the new edit.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("storage_medium")

$ws.Cells.Item(12, 1).Value = "Cryo-EM"
$ws.Cells.Item(12, 2).Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000333"

$ws.Cells.Item(13, 1).Value = "DMSO (serum)"
$ws.Cells.Item(13, 2).Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000125"

$ws.Cells.Item(14, 1).Value = "RNAlater"
$ws.Cells.Item(14, 2).Value = "http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C63348"

$meta = $wb.Worksheets.Item(".metadata")
$meta.Cells.Item(2, 3).Value = "2024-03-14T10:54:38-04:00"
